$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.878.45'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '2.083.79'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.27'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.395'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('E10').Value = '  +1.58%  '
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.83'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.26'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.777'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.36'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = '2.065.44'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '37.767.21'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.58'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('E20').Value = '  +3.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.33'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +2.02%  '
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.81'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.20'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.138'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.51'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +1.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.75'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.78'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.49%  '
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.52'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.45'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0987'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '99.39'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0218'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.76'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +7.21%  '
$ws.Range('D44').Value = '1.445.64'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.19'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.42'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').Value = '2.273.31'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.86'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.93%  '
